$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear out the old values in B2:B7, then set B2 to the single remaining value.
$ws.Range("B2:B7").ClearContents()
$ws.Range("B2").Value = "gh01"

# Move the active selection to B2 to match the saved view state.
$ws.Activate()
$ws.Range("B2").Select()
